$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 25; $r -le 38; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 + 1
}

$ws.Range("A2").Select()
